# Update cryptos list: apply latest price/volume data scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.111.47'
$ws.Range("E2").Value = '  +1.37%  '

$ws.Range("D3").Value = '2.630.36'
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.29%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.551'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("D9").Value = '2.629.11'
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.11%  '

$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("E12").Value = '  +0.94%  '

$ws.Range("E13").Value = '  -0.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.85%  '

$ws.Range("E15").Value = '  +4.55%  '

$ws.Range("D16").Value = '3.107.22'
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("D17").Value = '67.906.65'
$ws.Range("E17").Value = '  +1.09%  '

$ws.Range("D18").Value = '2.602.40'
$ws.Range("E18").Value = '  -0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '374.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.31%  '

$ws.Range("E20").Value = '  +2.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("E22").Value = '  -0.75%  '

$ws.Range("E23").Value = '  -1.05%  '

$ws.Range("E24").Value = '  -1.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.50%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("E28").Value = '  +3.17%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  -1.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '579.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.69%  '

$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.47%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.126'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '159.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.41%  '

$ws.Range("E39").Value = '  +0.19%  '

$ws.Range("E40").Value = '  +5.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.370'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.10%  '

$ws.Range("E43").Value = '  +4.36%  '

$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '0.0₆0323'
$ws.Range("E44").Value = '  +14.08%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.81%  '

$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.09%  '

$ws.Range("E49").Value = '  -0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.08'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.40%  '
